$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.642.67"
$ws.Range("E2").Value = "  -1.18%  "

# Row 3
$ws.Range("D3").Value = "2.192.58"
$ws.Range("E3").Value = "  -2.24%  "

# Row 4
$ws.Range("E4").Value = "  -0.29%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.11"
$ws.Range("E5").Value = "  +2.42%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.615"
$ws.Range("E6").Value = "  -0.76%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "75.33"
$ws.Range("E7").Value = "  -1.10%  "

# Row 8
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.585"
$ws.Range("E9").Value = "  -4.95%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.58"
$ws.Range("E10").Value = "  -1.76%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0916"
$ws.Range("E11").Value = "  -2.36%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.84"
$ws.Range("E12").Value = "  -2.04%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.102"
$ws.Range("E13").Value = "  +0.22%  "

# Row 14
$ws.Range("D14").Value = "2.519.47"
$ws.Range("E14").Value = "  -1.61%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.24"
$ws.Range("E15").Value = "  -3.33%  "

# Row 16
$ws.Range("D16").Value = "2.180.65"
$ws.Range("E16").Value = "  -2.65%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.774"
$ws.Range("E17").Value = "  -4.75%  "

# Row 18
$ws.Range("D18").Value = "42.550.28"
$ws.Range("E18").Value = "  -1.16%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000102"
$ws.Range("E19").Value = "  -2.84%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.88"
$ws.Range("E20").Value = "  -0.39%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.90"
$ws.Range("E21").Value = "  -1.73%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.48"
$ws.Range("E22").Value = "  -1.48%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.16"
$ws.Range("E23").Value = "  -2.89%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.49"
$ws.Range("E24").Value = "  -9.38%  "

# Row 25
$ws.Range("E25").Value = "  +0.03%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.49"
$ws.Range("E26").Value = "  -4.41%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.37"
$ws.Range("E27").Value = "  -0.13%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "38.74"
$ws.Range("E28").Value = "  +3.02%  "

# Row 29
$ws.Range("E29").Value = "  -0.81%  "

# Row 30
$ws.Range("E30").Value = "  -4.19%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.05"
$ws.Range("E31").Value = "  -0.62%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.12"
$ws.Range("E32").Value = "  -1.19%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0823"
$ws.Range("E33").Value = "  +3.49%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.16"
$ws.Range("E34").Value = "  -4.53%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.121"
$ws.Range("E35").Value = "  -1.50%  "

# Row 36
$ws.Range("E36").Value = "  -3.01%  "

# Row 37
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0338"
$ws.Range("E37").Value = "  +1.51%  "

# Row 38
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.23"
$ws.Range("E38").Value = "  -3.10%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.16"
$ws.Range("E39").Value = "  -7.64%  "

# Row 40
$ws.Range("E40").Value = "  -3.26%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.60"
$ws.Range("E41").Value = "  +12.41%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "59.29"
$ws.Range("E42").Value = "  -1.76%  "

# Row 43
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.195"
$ws.Range("E43").Value = "  -2.89%  "

# Row 44
$ws.Range("B44").Value = "THORChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.18"
$ws.Range("E44").Value = "  -7.38%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.90"
$ws.Range("E45").Value = "  -3.34%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0976"
$ws.Range("E46").Value = "  -1.66%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.24"
$ws.Range("E47").Value = "  -4.57%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.463"
$ws.Range("E48").Value = "  +4.00%  "

# Row 49
$ws.Range("E49").Value = "  -1.26%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.13"
$ws.Range("E50").Value = "  -2.09%  "

# Row 51
$ws.Range("E51").Value = "  -0.88%  "
